# Auto-generated edit script: updates cryptos list values per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.733.54'
$ws.Range('E2').Value = '  -0.91%  '
$ws.Range('D3').Value = '1.887.17'
$ws.Range('E3').Value = '  -0.78%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.7782'
$ws.Range('E5').Value = '  -6.17%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '241.43'
$ws.Range('E6').Value = '  -0.11%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.000'
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3156'
$ws.Range('E8').Value = '  -2.12%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '25.22'
$ws.Range('E9').Value = '  -5.79%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.06982'
$ws.Range('E10').Value = '  -0.57%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.08032'
$ws.Range('E11').Value = '  +0.06%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.7621'
$ws.Range('E12').Value = '  +1.54%  '
$ws.Range('D13').Value = '1.899.87'
$ws.Range('E13').Value = '  -0.11%  '
$ws.Range('E14').Value = '  +0.78%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '91.85'
$ws.Range('E15').Value = '  -0.86%  '
$ws.Range('D16').Value = '29.742.65'
$ws.Range('E16').Value = '  -0.91%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '13.76'
$ws.Range('E17').Value = '  -2.64%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '5.881'
$ws.Range('E18').Value = '  -0.96%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '241.76'
$ws.Range('E19').Value = '  -1.06%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.000007658'
$ws.Range('E20').Value = '  -1.49%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '1.001'
$ws.Range('E21').Value = '  +0.02%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '8.156'
$ws.Range('E22').Value = '  +16.98%  '
$ws.Range('D23').Value = '2.142.39'
$ws.Range('E23').Value = '  -0.62%  '
$ws.Range('E24').Value = '  -0.11%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.1631'
$ws.Range('E25').Value = '  +2.50%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.260'
$ws.Range('E26').Value = '  +0.51%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '165.19'
$ws.Range('E27').Value = '  -1.81%  '
$ws.Range('E28').Value = '  -1.89%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.036'
$ws.Range('E29').Value = '  -2.74%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.400'
$ws.Range('E30').Value = '  +1.99%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.529'
$ws.Range('E31').Value = '  +1.00%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.365'
$ws.Range('E32').Value = '  +2.07%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.05592'
$ws.Range('E33').Value = '  -2.62%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.016'
$ws.Range('E34').Value = '  -1.73%  '
$ws.Range('E35').Value = '  -2.49%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.7304'
$ws.Range('E36').Value = '  -0.68%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.004'
$ws.Range('E37').Value = '  +0.46%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.641'
$ws.Range('E38').Value = '  -3.06%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.01895'
$ws.Range('E39').Value = '  -0.55%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.766'
$ws.Range('E40').Value = '  -1.00%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.4376'
$ws.Range('E41').Value = '  -0.85%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '71.99'
$ws.Range('E42').Value = '  -0.28%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.786'
$ws.Range('E43').Value = '  -2.84%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.000'
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.8359'
$ws.Range('E45').Value = '  -1.03%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '101.77'
$ws.Range('E46').Value = '  +0.42%  '
$ws.Range('D47').Value = '1.015.23'
$ws.Range('E47').Value = '  +2.48%  '
$ws.Range('E48').Value = '  -2.51%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '9.854'
$ws.Range('E49').Value = '  +1.17%  '
$ws.Range('B50').Value = 'Aptos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.363'
$ws.Range('E50').Value = '  -3.38%  '
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = '2.037.68'
$ws.Range('E51').Value = '  -0.91%  '
